$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$needle = "You might also like "

for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null -and $val.Contains($needle)) {
        $cell.Value2 = $val.Replace($needle, "")
    }
}
